$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text type for numeric-looking strings in column D,
# then restore default "Normal" style so no stray number-format style lingers.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.307.19"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.875.08"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.06%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.7122"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "242.34"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.04%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.08035"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +3.41%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3147"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "24.98"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08218"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.95%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.883.55"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "94.87"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.12%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.248"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.7119"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.410"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +5.99%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000008518"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.99%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "29.307.88"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "243.84"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "2.129.01"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +0.28%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1559"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.035"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "162.27"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.51"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.501"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.56%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.399"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.299"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.05367"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.177"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -8.77%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.938"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.7639"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.87%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.177"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.690"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01874"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.15%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.265.35"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +3.04%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.752"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.80%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "6.442"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.98%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9152"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +3.44%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "112.36"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.92%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "73.92"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.99%  "

$ws.Range("E45").Value = "  +9.63%  "

$ws.Range("E46").Value = "  +0.00%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.026.33"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.5224"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("E49").Value = "  +0.24%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "9.477"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.4350"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
